$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "20190016830@my.xu.edu.ph"
$ws.Range("B2").Value = "James Jilhaney"
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:20190016830@my.xu.edu.ph", "", "", "20190016830@my.xu.edu.ph")
$ws.Range("A2").Font.Underline = $false
$ws.Range("A2").Font.ColorIndex = 0
"done"
